$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 364, pushing the existing row 364 (and
# everything below it) down by one. This grows the used range from
# A1:R459 to A1:R460, matching the new "dimension" in the edited file.
$ws.Rows.Item(364).Insert()

# Populate the newly inserted row 364 with the new weekly price record
# (Pepino ensalada, Región de Arica y Parinacota, fecha 2022-07-12).
$ws.Cells.Item(364, 1).Value = 6
$ws.Cells.Item(364, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(364, 3).Value = "Metropolitana"
$ws.Cells.Item(364, 4).Value = 44754
$ws.Cells.Item(364, 5).Value = 13
$ws.Cells.Item(364, 6).Value = 100112043
$ws.Cells.Item(364, 7).Value = "Pepino ensalada"
$ws.Cells.Item(364, 8).Value = "Sin especificar"
$ws.Cells.Item(364, 9).Value = "Primera"
$ws.Cells.Item(364, 10).Value = 400
$ws.Cells.Item(364, 11).Value = 17000
$ws.Cells.Item(364, 12).Value = 19000
$ws.Cells.Item(364, 13).Value = 17850
$ws.Cells.Item(364, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(364, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(364, 16).Value = 298
$ws.Cells.Item(364, 17).Value = 60
$ws.Cells.Item(364, 18).Value = "Hortaliza"
